$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---------
# The Overview sheet mirrors the per-locale status in columns E (zh-cn) and
# F (de-de); each locale sheet also carries its own "Status" column (C).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the status columns (regenerated report column sizing) --------
# ColumnWidth is rounded to whole pixels (width*6 px, plus the 5px cell
# padding) before it is stored, so the requested value is nudged down by
# that padding up front to land on the closest achievable stored width.
$targetWidth = 13.4101845877511
$inputWidth = $targetWidth - (5 / 6)

$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = $inputWidth
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $inputWidth
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $inputWidth
